$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 152 (weekly update adds a new
# price record for "Vega Modelo de Temuco" / "Bruselas (repollito)"),
# pushing the existing rows 152-157 down to 153-158.
$ws.Rows(152).Insert()

# Populate the newly inserted row 152 with this week's record.
$ws.Cells.Item(152, 1).Value = 10
$ws.Cells.Item(152, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value = "La Araucanía"
$ws.Cells.Item(152, 4).Value = 45041
$ws.Cells.Item(152, 5).Value = 9
$ws.Cells.Item(152, 6).Value = 100112035
$ws.Cells.Item(152, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 40
$ws.Cells.Item(152, 11).Value = 28000
$ws.Cells.Item(152, 12).Value = 28000
$ws.Cells.Item(152, 13).Value = 28000
$ws.Cells.Item(152, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(152, 15).Value = "Región Metropolitana"
$ws.Cells.Item(152, 16).Value = 1867
$ws.Cells.Item(152, 17).Value = 15
$ws.Cells.Item(152, 18).Value = "Hortaliza"

# Match the date-formatted style used by the other rows' "Fecha" column.
$ws.Cells.Item(152, 4).NumberFormat = $ws.Cells.Item(153, 4).NumberFormat
